$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
  2  = @{ D = "26.878.32";      E = "  +1.82%  " }
  3  = @{ D = "1.731.04";       E = "  +0.50%  " }
  4  = @{ D = "0.9963";         E = "  -0.13%  " }
  5  = @{ D = "241.14";         E = "  -0.63%  " }
  6  = @{ D = "0.9977";         E = "  -0.10%  " }
  7  = @{ D = "0.4833";         E = "  -1.09%  " }
  8  = @{ D = "0.2597";         E = "  -0.31%  " }
  9  = @{ D = "0.06188";        E = "  +0.02%  " }
  10 = @{ D = "1.722.81";       E = "  -0.02%  " }
  11 = @{ D = "15.93";          E = "  +2.50%  " }
  12 = @{ D = "0.06858";        E = "  -1.78%  " }
  13 = @{ D = "0.6058";         E = "  +1.14%  " }
  14 = @{ D = "4.475";          E = "  -0.94%  " }
  15 = @{ D = "77.08";          E = "  +0.02%  " }
  16 = @{ D = "0.9976" }
  17 = @{ D = "26.637.26";      E = "  +0.92%  " }
  18 = @{ D = "0.9959";         E = "  -0.19%  " }
  19 = @{ D = "0.000007147";    E = "  +0.17%  " }
  20 = @{ D = "11.40";          E = "  +1.00%  " }
  21 = @{ D = "1.943.43";       E = "  -0.28%  " }
  22 = @{ D = "4.424";          E = "  -0.43%  " }
  23 = @{ D = "8.493";          E = "  -0.10%  " }
  24 = @{ D = "5.081";          E = "  -0.21%  " }
  25 = @{ D = "140.31";         E = "  +1.74%  " }
  26 = @{ D = "15.29";          E = "  +0.32%  " }
  27 = @{ D = "1.794";          E = "  +3.36%  " }
  28 = @{ D = "106.86";         E = "  +0.51%  " }
  29 = @{ D = "1.370";          E = "  -2.36%  " }
  30 = @{ D = "3.983";          E = "  +2.05%  " }
  31 = @{ D = "0.07917";        E = "  -1.46%  " }
  32 = @{ D = "3.680";          E = "  +0.81%  " }
  33 = @{ D = "0.04546";        E = "  +1.35%  " }
  34 = @{ D = "2.588";          E = "  -0.61%  " }
  35 = @{ D = "1.004";          E = "  +0.73%  " }
  36 = @{ D = "0.6197";         E = "  -0.53%  " }
  37 = @{ D = "0.9281";         E = "  +0.62%  " }
  38 = @{ D = "2.020";          E = "  +2.89%  " }
  39 = @{ D = "2.449";          E = "  +2.59%  " }
  40 = @{ D = "0.9970";         E = "  -0.13%  " }
  41 = @{ D = "0.01498";        E = "  +1.37%  " }
  42 = @{ D = "5.669";          E = "  +4.83%  " }
  43 = @{ D = "99.77";          E = "  -0.13%  " }
  44 = @{ D = "0.3841";         E = "  +0.01%  " }
  45 = @{ D = "6.824";          E = "  -1.34%  " }
  46 = @{ D = "0.1158";         E = "  -0.28%  " }
  47 = @{ D = "0.05361";        E = "  -0.11%  " }
  48 = @{ D = "7.929";          E = "  +2.94%  " }
  49 = @{ D = "30.10";          E = "  -1.04%  " }
  50 = @{ D = "1.247";          E = "  +2.40%  " }
  51 = @{ D = "51.67";          E = "  +1.01%  " }
}

foreach ($row in $data.Keys) {
  $vals = $data[$row]

  $dCell = $ws.Range("D$row")
  $dCell.Value = "'" + $vals.D
  $dCell.ClearFormats()

  if ($vals.ContainsKey("E")) {
    $eCell = $ws.Range("E$row")
    $eCell.Value = "'" + $vals.E
    $eCell.ClearFormats()
  }
}
